$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 287, shifting existing rows 287-387 down to 289-389
$ws.Rows.Item(287).Resize(2).Insert()

# Row 287: new record
$ws.Cells.Item(287, 1).Value = 4
$ws.Cells.Item(287, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(287, 3).Value = "Los Lagos"
$ws.Cells.Item(287, 4).Value = 44795
$ws.Cells.Item(287, 5).Value = 10
$ws.Cells.Item(287, 6).Value = 100112008
$ws.Cells.Item(287, 7).Value = "Coliflor"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 500
$ws.Cells.Item(287, 11).Value = 2000
$ws.Cells.Item(287, 12).Value = 2000
$ws.Cells.Item(287, 13).Value = 2000
$ws.Cells.Item(287, 14).Value = "$/unidad"
$ws.Cells.Item(287, 15).Value = "Región del Maule"
$ws.Cells.Item(287, 16).Value = 2000
$ws.Cells.Item(287, 17).Value = 1
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# Row 288: new record
$ws.Cells.Item(288, 1).Value = 4
$ws.Cells.Item(288, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(288, 3).Value = "Los Lagos"
$ws.Cells.Item(288, 4).Value = 44795
$ws.Cells.Item(288, 5).Value = 10
$ws.Cells.Item(288, 6).Value = 100112008
$ws.Cells.Item(288, 7).Value = "Coliflor"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Segunda"
$ws.Cells.Item(288, 10).Value = 250
$ws.Cells.Item(288, 11).Value = 1700
$ws.Cells.Item(288, 12).Value = 1700
$ws.Cells.Item(288, 13).Value = 1700
$ws.Cells.Item(288, 14).Value = "$/unidad"
$ws.Cells.Item(288, 15).Value = "Región del Maule"
$ws.Cells.Item(288, 16).Value = 1700
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
